# sprint updated day 5
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day 4 summary - Total testcase Written
$ws.Range("C17").Value = 78

# Day 5 summary - Total testcase Written
$ws.Range("C25").Value = 133

# Move the active selection to C25 (Day 5 written-count cell)
$ws.Range("C25").Select()
